# Auto-generated edit script: refresh market-price / profit derived values
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 244.66667
$ws.Cells.Item(2, 9).Value = 147.14285
$ws.Cells.Item(2, 10).Value = 330
$ws.Cells.Item(2, 11).Value = 147.14285
$ws.Cells.Item(2, 12).Value = 330
$ws.Cells.Item(2, 13).Value = -34.14285000000001
$ws.Cells.Item(2, 14).Value = -556

$ws.Cells.Item(9, 8).Value = 89.76922999999999
$ws.Cells.Item(9, 9).Value = 88.916664
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 88.916664
$ws.Cells.Item(9, 12).Value = 100
$ws.Cells.Item(9, 13).Value = 80.083336
$ws.Cells.Item(9, 14).Value = -438

$ws.Cells.Item(15, 8).Value = 4048.7078
$ws.Cells.Item(15, 9).Value = 4048.7078
$ws.Cells.Item(15, 11).Value = 12146.1234
$ws.Cells.Item(15, 13).Value = -11977.1234

$ws.Cells.Item(113, 8).Value = 100004500
$ws.Cells.Item(113, 9).Value = 200000000
$ws.Cells.Item(113, 10).Value = 8999
$ws.Cells.Item(113, 11).Value = 200000000
$ws.Cells.Item(113, 12).Value = 8999
$ws.Cells.Item(113, 13).Value = -199996746
$ws.Cells.Item(113, 14).Value = -15507

$ws.Cells.Item(132, 8).Value = 12348855
$ws.Cells.Item(132, 9).Value = 13336363
$ws.Cells.Item(132, 11).Value = 40009089
$ws.Cells.Item(132, 13).Value = -40006559

$ws.Cells.Item(138, 8).Value = 2010.3776
$ws.Cells.Item(138, 10).Value = 2137.1704
$ws.Cells.Item(138, 12).Value = 6411.5112
$ws.Cells.Item(138, 14).Value = -16691.5112

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2902.516
$ws.Cells.Item(32, 9).Value = 3238.1667
$ws.Cells.Item(32, 11).Value = 3238.1667
$ws.Cells.Item(32, 13).Value = -2951.1667

$ws.Cells.Item(132, 8).Value = 2215.3408
$ws.Cells.Item(132, 9).Value = 1920.6216
$ws.Cells.Item(132, 11).Value = 5761.864799999999
$ws.Cells.Item(132, 13).Value = -3231.864799999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 100002210
$ws.Cells.Item(105, 9).Value = 111113336
$ws.Cells.Item(105, 10).Value = 2100
$ws.Cells.Item(105, 11).Value = 111113336
$ws.Cells.Item(105, 12).Value = 2100
$ws.Cells.Item(105, 13).Value = -111111589
$ws.Cells.Item(105, 14).Value = -5594

$ws.Cells.Item(134, 8).Value = 11120.954
$ws.Cells.Item(134, 9).Value = 7666
$ws.Cells.Item(134, 10).Value = 26668.25
$ws.Cells.Item(134, 11).Value = 22998
$ws.Cells.Item(134, 12).Value = 80004.75
$ws.Cells.Item(134, 13).Value = -20463
$ws.Cells.Item(134, 14).Value = -85074.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(25, 8).Value = 511
$ws.Cells.Item(25, 9).Value = 511
$ws.Cells.Item(25, 11).Value = 511
$ws.Cells.Item(25, 13).Value = -337

$ws.Cells.Item(31, 8).Value = 1448.3
$ws.Cells.Item(31, 9).Value = 1195
$ws.Cells.Item(31, 10).Value = 1757.8889
$ws.Cells.Item(31, 11).Value = 1195
$ws.Cells.Item(31, 12).Value = 1757.8889
$ws.Cells.Item(31, 13).Value = -900
$ws.Cells.Item(31, 14).Value = -2347.8889

$ws.Cells.Item(34, 8).Value = 1448.3
$ws.Cells.Item(34, 9).Value = 1195
$ws.Cells.Item(34, 10).Value = 1757.8889
$ws.Cells.Item(34, 11).Value = 1195
$ws.Cells.Item(34, 12).Value = 1757.8889
$ws.Cells.Item(34, 13).Value = -993
$ws.Cells.Item(34, 14).Value = -2161.8889

$ws.Cells.Item(39, 8).Value = 749.5
$ws.Cells.Item(39, 9).Value = 749.5
$ws.Cells.Item(39, 11).Value = 749.5
$ws.Cells.Item(39, 13).Value = -358.5

$ws.Cells.Item(49, 8).Value = 749.5
$ws.Cells.Item(49, 9).Value = 749.5
$ws.Cells.Item(49, 11).Value = 749.5
$ws.Cells.Item(49, 13).Value = -567.5

$ws.Cells.Item(112, 8).Value = 38000
$ws.Cells.Item(112, 10).Value = 38000
$ws.Cells.Item(112, 12).Value = 38000
$ws.Cells.Item(112, 14).Value = -40954

$ws.Cells.Item(132, 8).Value = 5716.3794
$ws.Cells.Item(132, 9).Value = 6031.7915
$ws.Cells.Item(132, 10).Value = 4202.4
$ws.Cells.Item(132, 11).Value = 18095.3745
$ws.Cells.Item(132, 12).Value = 12607.2
$ws.Cells.Item(132, 13).Value = -15565.3745
$ws.Cells.Item(132, 14).Value = -17667.2

$ws.Cells.Item(134, 8).Value = 10418098
$ws.Cells.Item(134, 9).Value = 13334718
$ws.Cells.Item(134, 10).Value = 1601.1428
$ws.Cells.Item(134, 11).Value = 40004154
$ws.Cells.Item(134, 12).Value = 4803.428400000001
$ws.Cells.Item(134, 13).Value = -40001619
$ws.Cells.Item(134, 14).Value = -9873.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(26, 8).Value = 280.5
$ws.Cells.Item(26, 9).Value = 260
$ws.Cells.Item(26, 10).Value = 301
$ws.Cells.Item(26, 11).Value = 780
$ws.Cells.Item(26, 12).Value = 903
$ws.Cells.Item(26, 13).Value = -492
$ws.Cells.Item(26, 14).Value = -1479

$ws.Cells.Item(34, 8).Value = 1599
$ws.Cells.Item(34, 10).Value = 1822.2941
$ws.Cells.Item(34, 12).Value = 5466.8823
$ws.Cells.Item(34, 14).Value = -5634.8823

$ws.Cells.Item(68, 8).Value = 2330.1
$ws.Cells.Item(68, 9).Value = 900
$ws.Cells.Item(68, 10).Value = 2687.625
$ws.Cells.Item(68, 11).Value = 2700
$ws.Cells.Item(68, 12).Value = 8062.875
$ws.Cells.Item(68, 13).Value = -1889
$ws.Cells.Item(68, 14).Value = -9684.875

$ws.Cells.Item(71, 8).Value = 2330.1
$ws.Cells.Item(71, 9).Value = 900
$ws.Cells.Item(71, 10).Value = 2687.625
$ws.Cells.Item(71, 11).Value = 8100
$ws.Cells.Item(71, 12).Value = 24188.625
$ws.Cells.Item(71, 13).Value = -4044
$ws.Cells.Item(71, 14).Value = -32300.625

$ws.Cells.Item(81, 8).Value = 2693.8262
$ws.Cells.Item(81, 10).Value = 2935.8333
$ws.Cells.Item(81, 12).Value = 8807.499899999999
$ws.Cells.Item(81, 14).Value = -11053.4999

$ws.Cells.Item(84, 8).Value = 2693.8262
$ws.Cells.Item(84, 10).Value = 2935.8333
$ws.Cells.Item(84, 12).Value = 26422.4997
$ws.Cells.Item(84, 14).Value = -37654.4997

$ws.Cells.Item(129, 8).Value = 21931084
$ws.Cells.Item(129, 10).Value = 8335043.5
$ws.Cells.Item(129, 12).Value = 25005130.5
$ws.Cells.Item(129, 14).Value = -25015130.5

$ws.Cells.Item(131, 8).Value = 20001406
$ws.Cells.Item(131, 10).Value = 1527.8837
$ws.Cells.Item(131, 12).Value = 4583.6511
$ws.Cells.Item(131, 14).Value = -14663.6511

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2922.647
$ws.Cells.Item(122, 9).Value = 1880.5454
$ws.Cells.Item(122, 10).Value = 4833.1665
$ws.Cells.Item(122, 11).Value = 5641.6362
$ws.Cells.Item(122, 12).Value = 14499.4995
$ws.Cells.Item(122, 13).Value = -3191.6362
$ws.Cells.Item(122, 14).Value = -19399.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1875.5
$ws.Cells.Item(7, 9).Value = 1714.1428
$ws.Cells.Item(7, 11).Value = 1714.1428
$ws.Cells.Item(7, 13).Value = -1602.1428

$ws.Cells.Item(40, 8).Value = 3318.8
$ws.Cells.Item(40, 9).Value = 2863
$ws.Cells.Item(40, 11).Value = 2863
$ws.Cells.Item(40, 13).Value = -2727

$ws.Cells.Item(101, 8).Value = 12500
$ws.Cells.Item(101, 10).Value = 12500
$ws.Cells.Item(101, 12).Value = 12500
$ws.Cells.Item(101, 14).Value = -18990

$ws.Cells.Item(126, 8).Value = 1875.5
$ws.Cells.Item(126, 9).Value = 1714.1428
$ws.Cells.Item(126, 11).Value = 5142.428400000001
$ws.Cells.Item(126, 13).Value = -2672.428400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 2615.9
$ws.Cells.Item(96, 9).Value = 2584.875
$ws.Cells.Item(96, 10).Value = 2740
$ws.Cells.Item(96, 11).Value = 2584.875
$ws.Cells.Item(96, 12).Value = 2740
$ws.Cells.Item(96, 13).Value = -1211.875
$ws.Cells.Item(96, 14).Value = -5486

$ws.Cells.Item(103, 8).Value = 25475.5
$ws.Cells.Item(103, 10).Value = 25475.5
$ws.Cells.Item(103, 12).Value = 25475.5
$ws.Cells.Item(103, 14).Value = -27819.5

$ws.Cells.Item(122, 8).Value = 21669158
$ws.Cells.Item(122, 9).Value = 21669158
$ws.Cells.Item(122, 11).Value = 65007474
$ws.Cells.Item(122, 13).Value = -65005024
